$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo in B1: "Cuentass" -> "Cuenta"
$ws.Range("B1").Value = "Cuenta"

# Update data values in row 2
$ws.Range("A2").Value = 555
$ws.Range("C2").Value = 200

# Update the active cell selection to H12
$ws.Range("H12").Select()
